$wb = $excel.ActiveWorkbook

$wsTests = $wb.Worksheets.Item("Tests")
$wsResult = $wb.Worksheets.Item("Result")

# Add the new row (WorkflowFile = "Framework\KillAllProcesses.xaml", ExpectedResult/Status = "Success")
# to both the "Tests" sheet (row 10, cols A:B) and the "Result" sheet (row 10, cols A:B).
$wsTests.Range("A10").Value = "Framework\KillAllProcesses.xaml"
$wsTests.Range("B10").Value = "Success"

$wsResult.Range("A10").Value = "Framework\KillAllProcesses.xaml"
$wsResult.Range("B10").Value = "Success"

# Update selections / active sheet to match the saved view state:
# "Tests" becomes the active (selected) tab, with B20 selected.
# "Result" becomes inactive, with A10 selected.
[void]$wsResult.Range("A10").Select()

[void]$wsTests.Activate()
[void]$wsTests.Range("B20").Select()
